# Paso de usarse un array de niveles a usarse una tabla hash de niveles.
# Apply the edits to the "Nivel Facil" sheet (the active sheet in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nivel Facil")
$ws.Activate()

# Move the "p" label from D1 to C1.
$ws.Range("D1").ClearContents()
$ws.Range("C1").Value = "p"

# Move the "d" label from H1 to I1.
$ws.Range("H1").ClearContents()
$ws.Range("I1").Value = "d"

# Add new "cv" label at B4.
$ws.Range("B4").Value = "cv"

# Remove the "f" label that used to sit at G5 (C5 keeps its own "S" value).
$ws.Range("G5").ClearContents()

# Add new "cp" label at E7.
$ws.Range("E7").Value = "cp"

# Change the E9 label from "c" to "cv".
$ws.Range("E9").Value = "cv"

# Update the active cell selection to E3.
$ws.Range("E3").Select()
